# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  ("Office Theme" colour scheme)
#   ppt/theme/theme2.xml  ("Integral" colour scheme - the one actually
#                          wired to the slide master / presentation, so
#                          the one that is visibly in effect)
# The authored change swaps the two colour schemes, so the deck's live
# theme (theme2.xml) ends up with the plain "Office" palette instead of
# "Integral", while the font scheme and format scheme (fills/lines/
# effects) are left untouched - they were already byte-identical between
# the two parts. Apply that by rewriting each theme colour slot via the
# standard ThemeColorScheme.Colors(i).RGB COM idiom (colours are
# addressed 1-12 in dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# order). PowerPoint's .RGB values are packed 0xBBGGRR (blue in the high
# byte), the reverse of the usual 0xRRGGBB reading order.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
